$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab (was "o554F-HW45.xpc")
$ws.Name = "o554F"

# Add a new averaged-intensity row (row 16), matching the formatting of the
# preceding data row (row 15) so it picks up the same border/bold/alignment
# style used for column A in this table.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(16, 3).Value = 0.9975163270514803
$ws.Cells.Item(16, 4).Value = 0.9935294117647059
$ws.Cells.Item(16, 5).Value = 0.9976763069520145
$ws.Cells.Item(16, 6).Value = 0.9975163270514803
$ws.Cells.Item(16, 7).Value = 0.9964705882352941
$ws.Cells.Item(16, 8).Value = 1
$ws.Cells.Item(16, 9).Value = 0.9976470588235294
$ws.Cells.Item(16, 10).Value = 0.9935294117647059
$ws.Cells.Item(16, 11).Value = 0.9956028593583601
$ws.Cells.Item(16, 12).Value = 0.9965595932049203
$ws.Cells.Item(16, 13).Value = 0.9971399488045041
